# Driver Script file updated
#
# TestSuite sheet: the Runmode flag ("Y"/"N") for every test-case row except
# the Login_Verification row is switched off ("N"), and the description for
# the login row is reworded from "execute" to "execution". The current
# selection is left on B16 (matching where the author's cursor ended up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip Runmode to "N" for every module except the first (Login_Verification,
# row 2) keeps its "Y". Do this before the B2 text edit so the shared-string
# table gets the new "N" entry ahead of the reworded login description,
# matching the workbook's string order.
$ws.Range("C3:C15").Value = "N"

# Fix wording: "execute" -> "execution" for the login test-case description.
$ws.Range("B2").Value = "All type of login execution"

# Leave the active selection on B16.
$ws.Range("B16").Select()
